# Income statement update: insert Depreciation & Amortization, EBITDA,
# Interest Expense, Earnings Before Tax, Income Tax Expense and Net Income
# build-out rows, and refresh the Gross Profit / R&D / SG&A figures that
# changed as part of the adjustments functionality.
#
# All the figures in this sheet are stored as literal text (not numbers),
# so every numeric-looking value is written with a leading apostrophe to
# force text, and then the cell style is reset back to "Normal" so Excel's
# automatic quote-prefix / number-format styling doesn't leak into the
# saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    $ws.Range($cell).Value = "'" + $text
    $ws.Range($cell).Style = "Normal"
}

# Row 4: Gross Profit - values updated
Set-TextValue "C4" "1,400"
Set-TextValue "D4" "1,700"
Set-TextValue "E4" "1,945"
Set-TextValue "F4" "2,233"
Set-TextValue "G4" "2,574"

# Row 5: was "Gross Profit Subtotal" / gross_profit_subtotal -> now "Research & Development" / r_d
$ws.Range("A5").Value = "  Research & Development"
$ws.Range("B5").Value = "r_d"
Set-TextValue "C5" "-100"
Set-TextValue "D5" "-120"
Set-TextValue "E5" "-130"
Set-TextValue "F5" "-137"
Set-TextValue "G5" "-158"

# Row 6: was "Research & Development" / r_d -> now "Selling, General & Administrative" / sg_a
$ws.Range("A6").Value = "  Selling, General & Administrative"
$ws.Range("B6").Value = "sg_a"
Set-TextValue "C6" "-200"
Set-TextValue "D6" "-250"
Set-TextValue "E6" "-261"
Set-TextValue "F6" "-281"
Set-TextValue "G6" "-274"

# Row 7: was "Selling, General & Administrative" / sg_a -> now "Depreciation & Amortization" / depreciation_amortization
$ws.Range("A7").Value = "  Depreciation & Amortization"
$ws.Range("B7").Value = "depreciation_amortization"
Set-TextValue "C7" "-30"
Set-TextValue "D7" "-35"
Set-TextValue "E7" "-36"
Set-TextValue "F7" "-37"
Set-TextValue "G7" "-38"

# Row 8: Total Operating Expenses - now populated with values
Set-TextValue "C8" "-330"
Set-TextValue "D8" "-405"
Set-TextValue "E8" "-426"
Set-TextValue "F8" "-456"
Set-TextValue "G8" "-471"

# Row 9: new row - EBITDA / ebitda
$ws.Range("A9").Value = "  EBITDA"
$ws.Range("B9").Value = "ebitda"
Set-TextValue "C9" "1,100"
Set-TextValue "D9" "1,330"
Set-TextValue "E9" "1,555"
Set-TextValue "F9" "1,814"
Set-TextValue "G9" "2,141"

# Row 10: Operating Income (EBIT) / operating_income, shifted down from row 9
$ws.Range("A10").Value = "  Operating Income (EBIT)"
$ws.Range("B10").Value = "operating_income"
Set-TextValue "C10" "1,070"
Set-TextValue "D10" "1,295"
Set-TextValue "E10" "1,519"
Set-TextValue "F10" "1,777"
Set-TextValue "G10" "2,103"

# Row 11: new row - Interest Expense / interest_expense
$ws.Range("A11").Value = "  Interest Expense"
$ws.Range("B11").Value = "interest_expense"
Set-TextValue "C11" "-50"
Set-TextValue "D11" "-60"
Set-TextValue "E11" "-63"
Set-TextValue "F11" "-66"
Set-TextValue "G11" "-69"

# Row 12: new row - Earnings Before Tax / ebt
$ws.Range("A12").Value = "  Earnings Before Tax"
$ws.Range("B12").Value = "ebt"
Set-TextValue "C12" "1,020"
Set-TextValue "D12" "1,235"
Set-TextValue "E12" "1,456"
Set-TextValue "F12" "1,711"
Set-TextValue "G12" "2,034"

# Row 13: new row - Income Tax Expense / taxes
$ws.Range("A13").Value = "  Income Tax Expense"
$ws.Range("B13").Value = "taxes"
Set-TextValue "C13" "-75"
Set-TextValue "D13" "-90"
Set-TextValue "E13" "-108"
Set-TextValue "F13" "-130"
Set-TextValue "G13" "-156"

# Row 14: new row - Net Income / net_income
$ws.Range("A14").Value = "  Net Income"
$ws.Range("B14").Value = "net_income"
Set-TextValue "C14" "945"
Set-TextValue "D14" "1,145"
Set-TextValue "E14" "1,348"
Set-TextValue "F14" "1,582"
Set-TextValue "G14" "1,878"
